$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rerun no longer produces the "RC6" (column M) counter, so drop it entirely.
$ws.Range("M1:M13").Delete()

# New counts from the rerun, one 1-hour interval per row (18 rows total,
# 6 more than before since the count window got extended).
$newData = @(
    @(0,  45393.33333333334, 1, 27, 2, 6,  1, 0, 105, 9, 31, 5),
    @(1,  45393.34027777778, 1, 32, 0, 1,  1, 3, 100, 5, 32, 2),
    @(2,  45393.34722222222, 1, 32, 0, 11, 1, 0, 92,  6, 31, 2),
    @(3,  45393.35416666666, 0, 15, 4, 11, 2, 2, 116, 4, 29, 1),
    @(4,  45393.36111111111, 0, 19, 3, 5,  4, 1, 114, 5, 25, 0),
    @(5,  45393.36805555555, 2, 27, 0, 10, 1, 0, 92,  1, 26, 3),
    @(6,  45393.66666666666, 0, 22, 0, 3,  0, 1, 105, 0, 17, 5),
    @(7,  45393.67361111111, 2, 21, 0, 2,  2, 1, 63,  0, 13, 3),
    @(8,  45393.68055555555, 0, 27, 0, 3,  0, 5, 71,  0, 9,  3),
    @(9,  45393.6875,        0, 39, 0, 2,  1, 1, 75,  0, 11, 2),
    @(10, 45393.69444444445, 0, 13, 0, 2,  3, 1, 58,  0, 16, 1),
    @(11, 45393.70138888889, 2, 35, 0, 5,  0, 2, 75,  0, 17, 5),
    @(12, 45393.83333333334, 0, 10, 0, 5,  0, 0, 43,  0, 17, 0),
    @(13, 45393.84027777778, 0, 11, 0, 2,  0, 1, 49,  0, 13, 0),
    @(14, 45393.84722222222, 0, 13, 0, 3,  0, 0, 60,  0, 12, 3),
    @(15, 45393.85416666666, 0, 22, 0, 5,  0, 0, 58,  0, 10, 0),
    @(16, 45393.86111111111, 0, 15, 0, 8,  2, 0, 50,  0, 7,  0),
    @(17, 45393.86805555555, 1, 21, 0, 4,  0, 0, 52,  0, 9,  0)
)

# Grab "known good" formatted cells to clone formatting from for any brand
# new rows we create below the old A1:L13 extent.
$ws.Cells.Item(2, 1).Copy()
$aFormatSrc = $true
$ws.Cells.Item(2, 2).Copy()

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $i + 2
    $vals = $newData[$i]
    for ($c = 1; $c -le $vals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# Columns A & B carry dedicated number formats (plain integer counter /
# date-time). Make sure every row - including the 6 brand-new ones at the
# bottom - carries that same formatting.
$ws.Cells.Item(2, 1).Copy()
$ws.Range("A2:A19").PasteSpecial(-4122)
$ws.Cells.Item(2, 2).Copy()
$ws.Range("B2:B19").PasteSpecial(-4122)
